$d = $word.ActiveDocument

# Locate the specific paragraph that currently reads:
# "Review Weekly will review draft of Weekly Status Report 3"
# (found under the "Nathan Stewart will:" bullet list) and fix it to read
# "Review draft of Weekly Status Report 3" by removing "Weekly will review ".
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*Review Weekly will review draft of Weekly Status Report 3*") {
        $r.Find.Execute("Weekly will review draft", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "draft", 2)
        break
    }
}
